$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), copying the formatting of the
# neighboring "sum" header (G1) so it keeps the same bold/border/
# center style used by the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column (H2), unformatted like the
# other numeric cells in row 2.
$ws.Range("H2").Value = 0
